$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Single-value cell updates (row, new text)
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "486"
$t.Cell(5,1).Range.Text  = "0.00001"
$t.Cell(8,1).Range.Text  = "0.00008"
$t.Cell(9,1).Range.Text  = "0.00019"
$t.Cell(10,1).Range.Text = "0.00021"
$t.Cell(11,1).Range.Text = "0.00027"
$t.Cell(12,1).Range.Text = "0.07294"

# Collapse the trailing multi-run statistics rows down to a single summary value
$t.Cell(44,1).Range.Text = "99.94"
$t.Cell(45,1).Range.Text = "0.07"
$t.Cell(46,1).Range.Text = "121"
